$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing cells whose values changed (rows 8-15) ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

$ws.Range("C9").Value = 16

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

$ws.Range("C12").Value = 10

$ws.Range("D13").Value = 8

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# --- Append two new rows (16 and 17) for line7/line8 shifted to extr positions ---
# Row 16 - copy the A15 cell formatting (bold/border/center) for the new A16 cell
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# Row 17
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false

$excel.CutCopyMode = $false

# --- Insert the new shared strings line7 / line8 and retarget B8/B9 to them ---
# (B8 previously showed extr1, B9 previously showed extr2; the table rows for
#  extr1..extr8 shift down by two, with line7/line8 taking their old slots.)
$ws.Range("B8").Value = "line7"
$ws.Range("B9").Value = "line8"
$ws.Range("B10").Value = "extr1"
$ws.Range("B11").Value = "extr2"
$ws.Range("B12").Value = "extr3"
$ws.Range("B13").Value = "extr4"
$ws.Range("B14").Value = "extr5"
$ws.Range("B15").Value = "extr6"
